# Add a new "2022-Q3" sheet (placed between "总计" and "2021-Q4") with the
# quarter's fund-holding detail, and add the corresponding summary row to
# the "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet    = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q3" sheet -------------------------------------
# Duplicate the existing "2021-Q4" sheet - the copy is placed right after the
# original and keeps all of the original 2021-Q4 data/formatting, so it just
# stays as the "2021-Q4" sheet. The *original* sheet object is then renamed
# to "2022-Q3" and its data is overwritten with the new quarter's figures.
# This keeps formatting/styles intact for both sheets with no extra work.
$q4Sheet.Copy($null, $q4Sheet)
$q4Copy = $wb.Worksheets.Item("2021-Q4 (2)")

$q4Sheet.Name = "2022-Q3"
$q4Copy.Name  = "2021-Q4"

$q3Sheet = $q4Sheet

# Force the numeric-looking text columns (B, D:G) to be stored as text,
# matching the source data (fund code / fund size / position% / weight% /
# market value are all reported as text strings in this workbook).
$q3Sheet.Range("B2:B7").NumberFormat = "@"
$q3Sheet.Range("D2:G7").NumberFormat = "@"

# Rows 5:7 are new (the template only had 3 data rows) - copy the index
# column's formatting (bold / border / centered) from row 2 onto them.
$q3Sheet.Range("A2").Copy()
$q3Sheet.Range("A5:A7").PasteSpecial(-4122)

$q3Data = @(
    @("003396", "东方红优享红利混合",               "14.04", "60.52", "5.00", "0.7020", 3),
    @("000480", "东方红新动力灵活配置混合",           "12.63", "78.30", "5.00", "0.6315", 5),
    @("169103", "东方红睿轩三年定期开放灵活配置混合", "11.31", "70.03", "5.00", "0.5655", 3),
    @("001564", "东方红京东大数据灵活配置混合",       "8.84",  "73.95", "5.00", "0.4420", 5),
    @("015769", "天弘低碳经济混合A",                 "1.19",  "79.16", "3.31", "0.0394", 8),
    @("015770", "天弘低碳经济混合C",                 "1.03",  "79.16", "3.31", "0.0341", 8)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    $rec = $q3Data[$i]
    $q3Sheet.Cells.Item($row, 1).Value = $i
    $q3Sheet.Cells.Item($row, 2).Value = $rec[0]
    $q3Sheet.Cells.Item($row, 3).Value = $rec[1]
    $q3Sheet.Cells.Item($row, 4).Value = $rec[2]
    $q3Sheet.Cells.Item($row, 5).Value = $rec[3]
    $q3Sheet.Cells.Item($row, 6).Value = $rec[4]
    $q3Sheet.Cells.Item($row, 7).Value = $rec[5]
    $q3Sheet.Cells.Item($row, 8).Value = $rec[6]
}

# --- 2. Update the "总计" (totals) sheet ------------------------------------
# Push the existing 2021-Q4 total row down to row 3 (copying the index
# column's formatting onto the newly used row), then write the new 2022-Q3
# total into row 2.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(3, 3).Value = 3
$totalSheet.Cells.Item(3, 4).Value = 0.22

$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 2.41
